$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F8").Value = -2
$ws.Range("F10").Value = -10
$ws.Range("F11").Value = -1
$ws.Range("F13").Value = -2
